$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "['Netherlands', 'Soviet Union']"
$ws.Range("G5").Value = "['France', 'Sweden']"
$ws.Range("G8").Value = "['Denmark', 'Sweden']"
$ws.Range("G9").Value = "['France', 'Sweden']"
$ws.Range("G10").Value = "['Denmark', 'Sweden']"
$ws.Range("G13").Value = "['France', 'Bulgaria']"
$ws.Range("G16").Value = "['France', 'Bulgaria']"
$ws.Range("G19").Value = "['Germany', 'Czech Republic']"
$ws.Range("G21").Value = "['Germany', 'Czech Republic']"
$ws.Range("G22").Value = "['Croatia', 'Portugal']"
$ws.Range("G24").Value = "['Romania', 'Portugal']"
$ws.Range("G26").Value = "['Romania', 'Portugal']"
$ws.Range("G27").Value = "['Belgium', 'Italy']"
$ws.Range("G29").Value = "['FR Yugoslavia', 'Norway']"
$ws.Range("G31").Value = "['Netherlands', 'France']"
$ws.Range("G32").Value = "['Spain', 'Greece']"
$ws.Range("G33").Value = "['Greece', 'Portugal']"
$ws.Range("G37").Value = "['Denmark', 'Sweden']"
$ws.Range("G38").Value = "['Germany', 'Czech Republic']"
$ws.Range("G39").Value = "['Netherlands', 'Czech Republic']"
$ws.Range("G42").Value = "['Croatia', 'Germany']"
$ws.Range("G43").Value = "['Netherlands', 'Romania']"
$ws.Range("G44").Value = "['Netherlands', 'Italy']"
$ws.Range("G46").Value = "['Spain', 'Russia']"
$ws.Range("G54").Value = "['Spain', 'Italy']"
$ws.Range("G56").Value = "['France', 'Romania', 'Switzerland']"
$ws.Range("G57").Value = "['France', 'Albania', 'Switzerland']"
$ws.Range("G59").Value = "['Germany', 'Poland', 'Northern Ireland']"
$ws.Range("G60").Value = "['Croatia', 'Czech Republic', 'Spain']"
$ws.Range("G61").Value = "['Croatia', 'Turkey', 'Spain']"
$ws.Range("G62").Value = "['Sweden', 'Belgium', 'Italy']"
$ws.Range("G63").Value = "['Republic of Ireland', 'Belgium', 'Italy']"
$ws.Range("G64").Value = "['Hungary', 'Portugal', 'Iceland']"
$ws.Range("G65").Value = "['Italy', 'Switzerland', 'Wales']"
$ws.Range("G66").Value = "['Finland', 'Belgium', 'Russia']"
$ws.Range("G68").Value = "['Finland', 'Belgium', 'Russia']"
$ws.Range("G69").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("G70").Value = "['England', 'Czech Republic', 'Croatia']"
$ws.Range("G71").Value = "['Sweden', 'Spain', 'Slovakia']"
$ws.Range("G72").Value = "['Germany', 'France', 'Portugal']"
$ws.Range("G73").Value = "['France', 'Hungary', 'Portugal']"
$ws.Range("G74").Value = "['Germany', 'France', 'Portugal']"
$ws.Range("G75").Value = "['France', 'Hungary', 'Portugal']"
$ws.Range("G76").Value = "['Germany', 'France', 'Portugal']"
$ws.Range("G77").Value = "['Germany', 'Scotland', 'Switzerland']"
$ws.Range("G78").Value = "['Germany', 'Hungary', 'Switzerland']"
$ws.Range("G79").Value = "['Spain', 'Italy', 'Albania']"
$ws.Range("G80").Value = "['Croatia', 'Spain', 'Italy']"
$ws.Range("G82").Value = "['Netherlands', 'Austria', 'France']"
$ws.Range("G83").Value = "['Ukraine', 'Romania', 'Belgium']"
$ws.Range("G84").Value = "['Ukraine', 'Slovakia', 'Belgium']"
$ws.Range("G85").Value = "['Slovakia', 'Belgium', 'Romania']"
$ws.Range("G87").Value = "['Turkey', 'Georgia', 'Portugal']"